$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the random "code" values in column B (category labels in column A are unchanged)
$ws.Range("B1").Value  = "9M"
$ws.Range("B3").Value  = "3V"
$ws.Range("B5").Value  = "0I"
$ws.Range("B6").Value  = "SK"
$ws.Range("B7").Value  = "KN3"
$ws.Range("B8").Value  = "5A"
$ws.Range("B9").Value  = "GM"

# New row 12: Ad copy / XB
$ws.Range("A12").Value = "Ad copy"
$ws.Range("B12").Value = "XB"

# Widen column A
$ws.Columns.Item(1).ColumnWidth = 16.5

# Move the active selection to F7
$ws.Range("F7").Select()
